# Edit: Tue, May 19, 2020  7:05:15 AM
#
# 1) Re-style the three tables (slides 14, 15, 16) with the new built-in
#    table style {3226A4C6-16EE-4F9A-8097-E28491BFD415} (they all currently
#    use the local custom style {0E73C0A6-E07B-4D94-AC2C-B20C592ECFDC}).
# 2) Re-colour the deck's theme (theme1.xml, backing the slide master /
#    Design) from the "Integral" ("Red Violet") palette to the standard
#    "Office Theme" ("Office") palette.

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newTableStyle = "{3226A4C6-16EE-4F9A-8097-E28491BFD415}"
$tableSlideIndexes = @(14, 15, 16)

foreach ($slideIdx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($slideIdx)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Theme colours --------------------------------------------------
# Office Theme colour scheme, in ThemeColorScheme.Item() order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
